$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @(2, 4, "90.728.17", $false),
    @(2, 5, "  -0.78%  ", $false),
    @(3, 4, "3.157.13", $false),
    @(3, 5, "  +1.46%  ", $false),
    @(4, 5, "  -0.05%  ", $false),
    @(5, 4, "215.77", $true),
    @(5, 5, "  -1.56%  ", $false),
    @(6, 4, "625.40", $true),
    @(6, 5, "  +1.07%  ", $false),
    @(7, 5, "  +26.63%  ", $false),
    @(8, 4, "0.368", $true),
    @(9, 4, "0.999", $true),
    @(9, 5, "  -0.05%  ", $false),
    @(10, 4, "3.154.95", $false),
    @(10, 5, "  +1.45%  ", $false),
    @(11, 4, "0.753", $true),
    @(11, 5, "  +11.77%  ", $false),
    @(12, 5, "  +5.97%  ", $false),
    @(13, 5, "  +5.40%  ", $false),
    @(14, 5, "  -5.14%  ", $false),
    @(15, 4, "35.32", $true),
    @(15, 5, "  +6.32%  ", $false),
    @(16, 4, "90.563.68", $false),
    @(16, 5, "  -0.61%  ", $false),
    @(17, 4, "3.741.84", $false),
    @(17, 5, "  +1.89%  ", $false),
    @(18, 4, "3.172.04", $false),
    @(18, 5, "  +2.41%  ", $false),
    @(19, 5, "  +2.89%  ", $false),
    @(20, 4, "14.68", $true),
    @(20, 5, "  +5.79%  ", $false),
    @(21, 2, "BitcoinCash", $false),
    @(21, 3, "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", $false),
    @(21, 4, "469.09", $true),
    @(21, 5, "  +8.01%  ", $false),
    @(22, 2, "PEPE", $false),
    @(22, 3, "https://coinranking.com/coin/03WI8NQPF+pepe-pepe", $false),
    @(22, 4, "0.0000213", $true),
    @(22, 5, "  -6.50%  ", $false),
    @(23, 4, "9.18", $true),
    @(23, 5, "  +7.01%  ", $false),
    @(24, 4, "5.17", $true),
    @(24, 5, "  +0.43%  ", $false),
    @(25, 4, "96.97", $true),
    @(25, 5, "  +15.18%  ", $false),
    @(26, 4, "5.91", $true),
    @(26, 5, "  +5.13%  ", $false),
    @(27, 4, "12.35", $true),
    @(27, 5, "  +3.75%  ", $false),
    @(28, 4, "3.320.79", $false),
    @(28, 5, "  +2.13%  ", $false),
    @(30, 4, "0.222", $true),
    @(30, 5, "  +55.83%  ", $false),
    @(31, 5, "  -2.44%  ", $false),
    @(32, 4, "9.28", $true),
    @(32, 5, "  +6.32%  ", $false),
    @(33, 4, "0.999", $true),
    @(33, 5, "  -1.53%  ", $false),
    @(34, 4, "27.19", $true),
    @(34, 5, "  +17.90%  ", $false),
    @(35, 4, "519.42", $true),
    @(35, 5, "  +0.06%  ", $false),
    @(36, 5, "  +5.56%  ", $false),
    @(37, 2, "Fetch.AI", $false),
    @(37, 3, "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet", $false),
    @(37, 4, "1.35", $true),
    @(37, 5, "  +3.90%  ", $false),
    @(38, 5, "  +0.13%  ", $false),
    @(39, 2, "Kaspa", $false),
    @(39, 3, "https://coinranking.com/coin/V8GxkwWow+kaspa-kas", $false),
    @(39, 4, "0.145", $true),
    @(39, 5, "  +2.90%  ", $false),
    @(40, 4, "3.61", $true),
    @(40, 5, "  -7.55%  ", $false),
    @(41, 4, "0.0912", $true),
    @(41, 5, "  +26.22%  ", $false),
    @(42, 4, "0.430", $true),
    @(42, 5, "  +15.92%  ", $false),
    @(43, 5, "  -0.37%  ", $false),
    @(44, 5, "  -0.05%  ", $false),
    @(45, 2, "Stacks", $false),
    @(45, 3, "https://coinranking.com/coin/mMPrMcB7+stacks-stx", $false),
    @(45, 4, "1.99", $true),
    @(45, 5, "  +5.98%  ", $false),
    @(46, 2, "ARBITRUM", $false),
    @(46, 3, "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", $false),
    @(46, 4, "0.742", $true),
    @(46, 5, "  +21.12%  ", $false),
    @(48, 4, "4.80", $true),
    @(48, 5, "  +14.09%  ", $false),
    @(49, 4, "150.67", $true),
    @(49, 5, "  +5.58%  ", $false),
    @(50, 5, "  +9.84%  ", $false),
    @(51, 5, "  +2.93%  ", $false)
)

foreach ($item in $changes) {
    $r = $item[0]
    $c = $item[1]
    $v = $item[2]
    $needsText = $item[3]
    $cell = $ws.Cells.Item($r, $c)
    if ($needsText) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $v
}